# Update the "想去人数" (want-to-go count) figures on the "展览" and
# "全部类型" sheets to reflect the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 795
$wsExhibit.Range("F6").Value = 88
$wsExhibit.Range("F7").Value = 282
$wsExhibit.Range("F8").Value = 4006
$wsExhibit.Range("F10").Value = 4718
$wsExhibit.Range("F11").Value = 522
$wsExhibit.Range("F12").Value = 1190
$wsExhibit.Range("F13").Value = 80

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 795
$wsAll.Range("F6").Value = 88
$wsAll.Range("F8").Value = 282
$wsAll.Range("F9").Value = 4006
$wsAll.Range("F11").Value = 4718
$wsAll.Range("F12").Value = 522
$wsAll.Range("F13").Value = 1190
$wsAll.Range("F14").Value = 80
